$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the footer/metadata rows (165-169): sample size, source, author, date notes.
#    This also shrinks the sheet dimension from A1:D169 down to A1:D163 automatically.
$ws.Range("A165:A169").EntireRow.Delete()

# 2) Rename header row (A1:D1) to snake_case machine-readable column names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 3) Title-case the Spanish connector words (de/del/la/las/los/el/y) in state & municipality names.
$ws.Range("B9").Value = "Mazapa De Madero"
$ws.Range("B11").Value = "San Cristóbal De Las Casas"
$ws.Range("B16").Value = "San Francisco Del Oro"
$ws.Range("B20").Value = "San Juan De Sabinas"
$ws.Range("A22").Value = "Ciudad De México"
$ws.Range("A32").Value = "Estado De México"
$ws.Range("B41").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B49").Value = "Silao De La Victoria"
$ws.Range("B50").Value = "Valle De Santiago"
$ws.Range("B52").Value = "Acapulco De Juárez"
$ws.Range("B54").Value = "Atoyac De Álvarez"
$ws.Range("B55").Value = "Ayutla De Los Libres"
$ws.Range("B56").Value = "Chilapa De Álvarez"
$ws.Range("B57").Value = "Coyuca De Catalán"
$ws.Range("B59").Value = "Cutzamala De Pinzón"
$ws.Range("B64").Value = "Técpan De Galeana"
$ws.Range("B68").Value = "Huasca De Ocampo"
$ws.Range("B70").Value = "Santiago De Anaya"
$ws.Range("B72").Value = "Tula De Allende"
$ws.Range("B74").Value = "Ahualulco De Mercado"
$ws.Range("B92").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B100").Value = "Mier Y Noriega"
$ws.Range("B105").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B106").Value = "Ixtlán De Juárez"
$ws.Range("B108").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B111").Value = "San Antonino El Alto"
$ws.Range("B118").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B119").Value = "Tlacolula De Matamoros"
$ws.Range("B120").Value = "Villa De Chilapa De Díaz"
$ws.Range("B129").Value = "Cadereyta De Montes"
$ws.Range("B130").Value = "Pinal De Amoles"
$ws.Range("B133").Value = "San Juan Del Río"
$ws.Range("B136").Value = "Axtla De Terrazas"
$ws.Range("B157").Value = "Poza Rica De Hidalgo"

# 4) Refresh the percentage column (column D) with recomputed floating point values.
$ws.Range("D2").Value = 0.009433962264150945
$ws.Range("D3").Value = 0.009433962264150945
$ws.Range("D17").Value = 0.009433962264150945
$ws.Range("D24").Value = 0.009433962264150945
$ws.Range("D30").Value = 0.009433962264150945
$ws.Range("D31").Value = 0.009433962264150945
$ws.Range("D33").Value = 0.009433962264150945
$ws.Range("D35").Value = 0.009433962264150945
$ws.Range("D52").Value = 0.009433962264150945
$ws.Range("D59").Value = 0.009433962264150945
$ws.Range("D66").Value = 0.09905660377358493
$ws.Range("D69").Value = 0.009433962264150945
$ws.Range("D70").Value = 0.009433962264150945
$ws.Range("D78").Value = 0.009433962264150945
$ws.Range("D84").Value = 0.009433962264150945
$ws.Range("D94").Value = 0.009433962264150945
$ws.Range("D95").Value = 0.009433962264150945
$ws.Range("D101").Value = 0.009433962264150945
$ws.Range("D103").Value = 0.009433962264150945
$ws.Range("D114").Value = 0.009433962264150945
$ws.Range("D120").Value = 0.009433962264150945
$ws.Range("D122").Value = 0.009433962264150945
$ws.Range("D131").Value = 0.009433962264150945
$ws.Range("D133").Value = 0.009433962264150945
$ws.Range("D145").Value = 0.009433962264150945
$ws.Range("D150").Value = 0.009433962264150945
$ws.Range("D156").Value = 0.009433962264150945
